$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly price records appended to the "Fruta, Vega Monumental
# Concepción - Durazno" sheet (Royal Glory variety, week of 2022-12-16 /
# serial 44911), one row per quality grade.

$rows = @(
    @{ L = "Especial"; N = 15000; O = 15000; P = 15000; S = 1000 },
    @{ L = "Primera";  N = 14000; O = 14000; P = 14000; S = 933 },
    @{ L = "Segunda";  N = 12000; O = 12000; P = 12000; S = 800 }
)

$startRow = 216

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = 11
    $ws.Cells.Item($r, 2).Value = "Vega Monumental Concepción"
    $ws.Cells.Item($r, 3).Value = "Bíobío"
    $ws.Cells.Item($r, 4).Value = 44911
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = 8
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100103
    $ws.Cells.Item($r, 8).Value = "Frutos de hueso (carozo)"
    $ws.Cells.Item($r, 9).Value = 100103004
    $ws.Cells.Item($r, 10).Value = "Durazno"
    $ws.Cells.Item($r, 11).Value = "Royal Glory"
    $ws.Cells.Item($r, 12).Value = $data.L
    $ws.Cells.Item($r, 13).Value = 50
    $ws.Cells.Item($r, 14).Value = $data.N
    $ws.Cells.Item($r, 15).Value = $data.O
    $ws.Cells.Item($r, 16).Value = $data.P
    $ws.Cells.Item($r, 17).Value = "$/caja 15 kilos empedrada"
    $ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
    $ws.Cells.Item($r, 19).Value = $data.S
    $ws.Cells.Item($r, 20).Value = 15
}
